$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. LoginPage : just move the lingering selection from B5 to A5
# ---------------------------------------------------------------------------
$loginPage = $wb.Worksheets.Item("LoginPage")
[void]$loginPage.Range("A5").Select()

# ---------------------------------------------------------------------------
# 2. ManageProductsPage : replace the old 1-column "Title" sheet with the
#    full product-entry form (10 columns of headers + one sample row)
# ---------------------------------------------------------------------------
$products = $wb.Worksheets.Item("ManageProductsPage")

$products.Range("A1").Value = "Title"
$products.Range("B1").Value = "Tag"
$products.Range("C1").Value = "Weight Value"
$products.Range("D1").Value = "Maximum Quantity Can Order"
$products.Range("E1").Value = "Price"
$products.Range("F1").Value = "MRP"
$products.Range("G1").Value = "Stock Availability(Kg)"
$products.Range("H1").Value = "Purchase Price"
$products.Range("I1").Value = "Description"
$products.Range("J1").Value = "Image Location"

$products.Range("A2").Value = "Cake"
$products.Range("B2").Value = "Combo"
$products.Range("C2").Value = 10
$products.Range("D2").Value = 5
$products.Range("E2").Value = 1000
$products.Range("F2").Value = 950
$products.Range("G2").Value = 20
$products.Range("H2").Value = 1000
$products.Range("I2").Value = "Brand new cake additions"
$products.Range("J2").Value = "C:\\Users\\angit\\Downloads\\Desktop\\cake pic.jpg"

$products.Columns.Item(1).ColumnWidth = 14.109375
$products.Columns.Item(3).ColumnWidth = 11.77734375
$products.Columns.Item(4).ColumnWidth = 25.6640625
$products.Columns.Item(7).ColumnWidth = 17.88671875
$products.Columns.Item(8).ColumnWidth = 12.77734375
$products.Columns.Item(9).ColumnWidth = 22.109375
$products.Columns.Item(10).ColumnWidth = 42.77734375

$products.PageSetup.Orientation = 1

[void]$products.Range("J2").Select()

# ---------------------------------------------------------------------------
# 3. New sheets : AddOffercode, SearchOfferCode, EditOfferCode,
#    DeleteOfferCode, AddLocation (appended, in this order, after
#    ManageProductsPage)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$addOffer = $wb.Worksheets.Add($null, $lastSheet)
$addOffer.Name = "AddOffercode"
$lastSheet = $addOffer

$searchOffer = $wb.Worksheets.Add($null, $lastSheet)
$searchOffer.Name = "SearchOfferCode"
$lastSheet = $searchOffer

$editOffer = $wb.Worksheets.Add($null, $lastSheet)
$editOffer.Name = "EditOfferCode"
$lastSheet = $editOffer

$deleteOffer = $wb.Worksheets.Add($null, $lastSheet)
$deleteOffer.Name = "DeleteOfferCode"
$lastSheet = $deleteOffer

$addLocation = $wb.Worksheets.Add($null, $lastSheet)
$addLocation.Name = "AddLocation"
$lastSheet = $addLocation

# ---------------------------------------------------------------------------
# 3a. AddOffercode
# ---------------------------------------------------------------------------
$addOffer.Range("A1").Value = "Offer Code"
$addOffer.Range("B1").Value = "Percentage"
$addOffer.Range("C1").Value = "Amount"
$addOffer.Range("D1").Value = "Description"
$addOffer.Range("E1").Value = "Image"

$addOffer.Range("A2").Value = "T140"
$addOffer.Range("B2").Value = 10
$addOffer.Range("C2").Value = 500
$addOffer.Range("D2").Value = "Offer code T136 description"
$addOffer.Range("E2").Value = "C:\\Users\\angit\\Downloads\\Desktop\\cake pic.jpg"

$addOffer.Columns.Item(1).ColumnWidth = 10
$addOffer.Columns.Item(4).ColumnWidth = 24
$addOffer.Columns.Item(5).ColumnWidth = 45.109375

[void]$addOffer.Range("A2").Select()

# ---------------------------------------------------------------------------
# 3b. SearchOfferCode
# ---------------------------------------------------------------------------
$searchOffer.Range("A1").Value = "Offer Code"
$searchOffer.Range("A2").Value = "T138"

$searchOffer.Columns.Item(1).ColumnWidth = 10

[void]$searchOffer.Range("B6").Select()

# ---------------------------------------------------------------------------
# 3c. EditOfferCode
# ---------------------------------------------------------------------------
$editOffer.Range("A1").Value = "Offer Code"
$editOffer.Range("B1").Value = "Description"
$editOffer.Range("C1").Value = "Amount"

$editOffer.Range("A2").Value = "T138"
$editOffer.Range("B2").Value = "Edited description for Offer code T138"
$editOffer.Range("C2").Value = 1000

$editOffer.Columns.Item(2).ColumnWidth = 32.6640625

[void]$editOffer.Range("B10").Select()

# ---------------------------------------------------------------------------
# 3d. DeleteOfferCode
# ---------------------------------------------------------------------------
$deleteOffer.Range("A1").Value = "Offer Code"
$deleteOffer.Range("A2").Value = "T125"

[void]$deleteOffer.Range("M18").Select()

# ---------------------------------------------------------------------------
# 3e. AddLocation
# ---------------------------------------------------------------------------
$addLocation.Range("A1").Value = "Country"
$addLocation.Range("B1").Value = "State"
$addLocation.Range("C1").Value = "Location"
$addLocation.Range("D1").Value = "Delivery Charge"

$addLocation.Range("A2").Value = "United Kingdom"
$addLocation.Range("B2").Value = 3815
$addLocation.Range("C2").Value = "Lynch Street"
$addLocation.Range("D2").Value = 100

$addLocation.Columns.Item(1).ColumnWidth = 14
$addLocation.Columns.Item(3).ColumnWidth = 11

$addLocation.PageSetup.Orientation = 1

[void]$addLocation.Range("E5").Select()

# ---------------------------------------------------------------------------
# 4. Make AddLocation the active sheet/tab, matching the saved view state
# ---------------------------------------------------------------------------
[void]$addLocation.Select()
